# Add the new "Gameweeks" related columns to the Challenges import template:
#   S: "Show Statistics Continuously" (boolean-looking text flag)
#   T: "Gameweek" (numeric)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Challenges")

# Headers (row 1)
$ws.Range("S1").Value = "Show Statistics Continuously"
$ws.Range("T1").Value = "Gameweek"

# Sample data (row 2).
# The leading apostrophe forces the literal text "true" instead of Excel's
# auto-boolean conversion; reset the style afterwards so no quote-prefix
# formatting is left behind on the cell.
$ws.Range("S2").Value = "'true"
$ws.Range("S2").Style = "Normal"
$ws.Range("T2").Value = 1
